$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "variables_1080","variables_1081","variables_1082","variables_1083","variables_1084",
    "variables_1085","variables_1086","variables_1087","variables_1088","variables_1089",
    "variables_1090","variables_1091","variables_1092","variables_1093","variables_1094",
    "variables_1095","variables_1096","variables_1097","variables_1098","variables_1099",
    "variables_1100","variables_1101","variables_1102","variables_1103","variables_1104",
    "variables_1105","variables_1106","variables_1107","variables_1108","variables_1109",
    "variables_1110","variables_1111","variables_1112","variables_1113","variables_1114",
    "variables_1150","variables_1151","variables_1152","variables_1153","variables_1154",
    "variables_1155","variables_1156","variables_1157"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A4").Value = "misc_long_term"
}
